$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A = "14-06-2017"; B = "10:36"; C = "Knie extensie";             D = 0.0002199911226851852; E = 1 },
    @{ A = "14-06-2017"; B = "10:37"; C = "Leunen naar grond";         D = 0.0001773612615740741; E = 0 },
    @{ A = "14-06-2017"; B = "10:37"; C = "Naar voren leunen";         D = 0.00005244998842592592; E = 0 },
    @{ A = "14-06-2017"; B = "10:59"; C = "Knien optillen";            D = 0.01494191704861111;   E = 0 },
    @{ A = "14-06-2017"; B = "11:00"; C = "Knien en handen optillen";  D = 0.00005223804398148148; E = 0 },
    @{ A = "14-06-2017"; B = "11:00"; C = "Beide benen strekken";      D = 0.00005206115740740741; E = 0 }
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $r++
}

# Reuse the existing cell styles (s=3 for the duration column, s=4 for
# the trailing empty column) instead of letting NumberFormat create new
# style entries.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D4:D9").PasteSpecial(-4122) | Out-Null

$ws.Range("F2").Copy() | Out-Null
$ws.Range("F4:F9").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
